# Integrate parent document retriever
# The "Location" column (C) previously stored the short form "Hồ Chí Minh"
# for every data row. Replace it with the fuller "Thành phố Hồ Chí Minh"
# across all rows that currently hold that value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Hồ Chí Minh") {
        $cell.Value = "Thành phố Hồ Chí Minh"
    }
}
